$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.520.90"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "'1.904.36"
$ws.Range("E3").Value = "  -0.67%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'338.31"
$ws.Range("E5").Value = "  +4.09%  "

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("D7").Value = "'0.4767"
$ws.Range("E7").Value = "  -1.01%  "

# Row 8
$ws.Range("D8").Value = "'0.4002"
$ws.Range("E8").Value = "  -1.78%  "

# Row 9
$ws.Range("D9").Value = "'0.08039"
$ws.Range("E9").Value = "  -2.27%  "

# Row 10
$ws.Range("D10").Value = "'0.9904"
$ws.Range("E10").Value = "  -2.34%  "

# Row 11
$ws.Range("D11").Value = "'23.19"
$ws.Range("E11").Value = "  -0.92%  "

# Row 12
$ws.Range("D12").Value = "'1.902.43"
$ws.Range("E12").Value = "  -0.90%  "

# Row 13
$ws.Range("D13").Value = "'5.917"
$ws.Range("E13").Value = "  -2.65%  "

# Row 14
$ws.Range("D14").Value = "'7.105"
$ws.Range("E14").Value = "  -1.97%  "

# Row 15
$ws.Range("D15").Value = "'89.07"
$ws.Range("E15").Value = "  -2.54%  "

# Row 16
$ws.Range("D16").Value = "'0.06831"

# Row 17
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("E18").Value = "  -1.77%  "

# Row 19
$ws.Range("D19").Value = "'17.35"
$ws.Range("E19").Value = "  -1.53%  "

# Row 20
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21").Value = "'29.529.89"
$ws.Range("E21").Value = "  +0.27%  "

# Row 22
$ws.Range("D22").Value = "'5.508"

# Row 23
$ws.Range("D23").Value = "'11.58"
$ws.Range("E23").Value = "  -1.70%  "

# Row 24
$ws.Range("D24").Value = "'2.159"
$ws.Range("E24").Value = "  -0.77%  "

# Row 25
$ws.Range("D25").Value = "'2.137.35"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26
$ws.Range("D26").Value = "'157.04"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27
$ws.Range("D27").Value = "'6.483"
$ws.Range("E27").Value = "  -2.64%  "

# Row 28
$ws.Range("D28").Value = "'19.53"
$ws.Range("E28").Value = "  -2.48%  "

# Row 29
$ws.Range("D29").Value = "'2.049"
$ws.Range("E29").Value = "  -3.04%  "

# Row 30
$ws.Range("D30").Value = "'118.97"
$ws.Range("E30").Value = "  -1.29%  "

# Row 31
$ws.Range("D31").Value = "'0.9965"
$ws.Range("E31").Value = "  -2.17%  "

# Row 32
$ws.Range("D32").Value = "'0.09530"
$ws.Range("E32").Value = "  -0.70%  "

# Row 33
$ws.Range("D33").Value = "'5.463"
$ws.Range("E33").Value = "  -3.96%  "

# Row 34
$ws.Range("D34").Value = "'3.538"
$ws.Range("E34").Value = "  -0.36%  "

# Row 35
$ws.Range("D35").Value = "'1.386"
$ws.Range("E35").Value = "  +0.95%  "

# Row 36
$ws.Range("D36").Value = "'0.06461"
$ws.Range("E36").Value = "  +5.85%  "

# Row 37
$ws.Range("D37").Value = "'0.02239"
$ws.Range("E37").Value = "  -1.95%  "

# Row 38
$ws.Range("E38").Value = "  +0.94%  "

# Row 39
$ws.Range("E39").Value = "  -2.78%  "

# Row 40
$ws.Range("D40").Value = "'10.53"
$ws.Range("E40").Value = "  -2.77%  "

# Row 41
$ws.Range("D41").Value = "'7.759"
$ws.Range("E41").Value = "  -3.85%  "

# Row 42
$ws.Range("D42").Value = "'0.1819"
$ws.Range("E42").Value = "  -1.53%  "

# Row 43
$ws.Range("D43").Value = "'2.449"
$ws.Range("E43").Value = "  +1.62%  "

# Row 44
$ws.Range("D44").Value = "'1.268"
$ws.Range("E44").Value = "  -0.80%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07414"
$ws.Range("E45").Value = "  -2.42%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.13"
$ws.Range("E46").Value = "  -2.78%  "

# Row 47
$ws.Range("D47").Value = "'0.5479"
$ws.Range("E47").Value = "  -2.07%  "

# Row 48
$ws.Range("D48").Value = "'1.946"
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("E49").Value = "  -2.00%  "

# Row 50
$ws.Range("D50").Value = "'2.376"
$ws.Range("E50").Value = "  -2.19%  "

# Row 51
$ws.Range("D51").Value = "'71.10"
$ws.Range("E51").Value = "  -1.78%  "
